$d = $word.ActiveDocument

# Each tuple: Table index, Row index, Column index (all 1-based)
# Sets cell shading fill color to white (0xFFFFFF) to match updated figures.
$targets = @(
    @(1, 2, 4),
    @(1, 4, 4),
    @(1, 5, 3),
    @(1, 6, 3),
    @(1, 6, 4),
    @(1, 8, 3),
    @(1, 9, 2),
    @(1, 9, 3),
    @(1, 9, 4),
    @(1, 11, 3),
    @(1, 11, 4),
    @(1, 13, 2),
    @(1, 13, 3),
    @(1, 13, 4),
    @(2, 2, 2),
    @(2, 2, 4),
    @(2, 4, 2),
    @(2, 4, 4),
    @(2, 5, 2),
    @(2, 5, 4),
    @(2, 6, 3),
    @(2, 6, 4),
    @(2, 7, 4),
    @(2, 9, 3),
    @(2, 9, 4),
    @(2, 10, 2),
    @(2, 10, 4),
    @(2, 11, 3),
    @(2, 11, 4),
    @(2, 12, 4),
    @(2, 13, 3),
    @(2, 13, 4),
)

foreach ($tup in $targets) {
    $tblIdx = $tup[0]
    $rowIdx = $tup[1]
    $colIdx = $tup[2]
    $tbl = $d.Tables.Item($tblIdx)
    $cell = $tbl.Cell($rowIdx, $colIdx)
    $cell.Shading.BackgroundPatternColor = 0xFFFFFF
}

Write-Host "Done. Updated" $targets.Count "cells."
